$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty / duplicated "OBSERVACIONES" (L) column with
# unique placeholder values (UNO..CATORCE) so there are no more duplicated
# rows of test data (PRUEBAS / SADÑLNSA / GAM are replaced).
$ws.Range("L5").Value  = "UNO"
$ws.Range("L6").Value  = "DOS"
$ws.Range("L7").Value  = "TRES"
$ws.Range("L8").Value  = "CUATRO"
$ws.Range("L9").Value  = "CINCO"
$ws.Range("L10").Value = "SEIS"
$ws.Range("L11").Value = "SIETE"
$ws.Range("L12").Value = "OCHO"
$ws.Range("L13").Value = "NUEVE"
$ws.Range("L14").Value = "DIEZ"
$ws.Range("L15").Value = "ONCE"
$ws.Range("L16").Value = "DOCE"
$ws.Range("L17").Value = "TRECE"
$ws.Range("L18").Value = "CATORCE"

# K column (CAPATAZ) values are unchanged in content (DAVID / PEPE / SERGIO)
# but re-set explicitly to keep the foreman names consistent now that the
# duplicate test rows are gone.
$ws.Range("K5").Value  = "DAVID"
$ws.Range("K11").Value = "PEPE"
$ws.Range("K12").Value = "PEPE"

# Totals row: turn the per-column SUM formulas into one shared formula
# (fill C19 across to I19 with a single relative formula).
$ws.Range("C19:I19").Formula = "=SUM(C5:C18)"

# Restore J36's own formula (kept intact, just re-applied so it survives the
# shared-formula renumbering above).
$ws.Range("J36").Formula = "=F36*H36"

# Move the active selection to L18, matching where the work left off.
$ws.Range("L18").Select() | Out-Null
